$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before current row 17 (old rows 17-19 become 19-21,
# old rows 20-48 become 22-50). We'll then overwrite rows 17-21 with the
# updated / reorganised data for the endogenous croc/turtle/komodo rows.
$ws.Rows.Item(17).Resize(2).EntireRow.Insert()

# Row 17: eCRHBV (was eCRHBV1, renamed/merged)
$ws.Range("A17").Value = "eCRHBV"
$ws.Range("B17").Value = "eCRHBV-con"
$ws.Range("C17").Value = "Endogenous crocodile hepatitis B virus "
$ws.Range("D17").Value = "Herpetohepadnavirus"
$ws.Range("E17").Value = "Crocodylus"
$ws.Range("F17").Value = "crocodiles"

# Row 18: eCRHBV2 (unchanged content, new position)
$ws.Range("A18").Value = "eCRHBV2"
$ws.Range("B18").Value = "eCRHBV2-con"
$ws.Range("C18").Value = "Endogenous crocodile hepatitis B virus 2"
$ws.Range("D18").Value = "Herpetohepadnavirus"
$ws.Range("E18").Value = "Crocodylus"
$ws.Range("F18").Value = "crocodiles and gharial"

# Match the font formatting of column F on these two rows to the rest of
# column E/F in this style family (style "s=2", not the heading-row style
# that Insert() propagated from row 16 above).
$ws.Range("E17").Copy() | Out-Null
$ws.Range("F17").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Value = "crocodiles"
$ws.Range("F18").Value = "crocodiles and gharial"

# Row 19: eDRHPV (new entry - endogenous komodo dragon hepatitis B virus)
$ws.Range("A19").Value = "eDRHPV"
$ws.Range("B19").Value = "eDRHPV-con"
$ws.Range("C19").Value = "Endogenous comodo dragon hepatitis B virus"
$ws.Range("D19").Value = "Herpetohepadnavirus"
$ws.Range("E19").Value = "Comodo dragon"
$ws.Range("F19").Value = ""

# Row 20: eTHBV (was eTHBV1, renamed)
$ws.Range("A20").Value = "eTHBV"
$ws.Range("B20").Value = "eTHBV-con"
$ws.Range("C20").Value = "Endogenous turtle hepatitis B virus"
$ws.Range("D20").Value = "Herpetohepadnavirus"
$ws.Range("E20").Value = "Chrysemys"
$ws.Range("F20").Value = "turtles"

# Row 21: left blank as a separator row
$ws.Range("A21:F21").Value = ""

# Update defined name for the filter database to the new extent
$wb.Names.Item("_xlnm._FilterDatabase").RefersToR1C1 = "=Sheet1!R1C1:R50C6"

# Update the autofilter range to cover the new data extent
$ws.Range("A1:F50").AutoFilter()

# Adjust the view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A21:XFD21").Select()
